$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in column B
$ws.Range("B2").Value = 29.5606346924565
$ws.Range("B3").Value = -95.088395130234403

# Clear column C cells entirely (no shifting of later columns)
$ws.Range("C1:C17").Clear()

# Update the active selection
$ws.Range("B20").Select()
